$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clarify the header label: "V/N" -> "V/U"
$ws.Range("D3").Value = "V/U"

# Add a small legend in column F explaining the V/U codes
$ws.Range("F3").Value = "V = Vanilla"
$ws.Range("F4").Value = "U = User-created"

# Bold the legend text on F4 to match the rest of the header row styling
$ws.Range("F4").Font.Bold = $true

# Widen column F so the legend text isn't truncated
$ws.Range("F1").EntireColumn.ColumnWidth = 16.666666666666668

# Leave the selection on the newly-edited cell, as the author did
$ws.Range("F3").Select() | Out-Null
